$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round Q7 / R7 to whole numbers
$ws.Range("Q7").Value = 527229
$ws.Range("R7").Value = 6908169

# Delete contents of Z7 and AB7 (Starttid / Sluttid for row 7)
$ws.Range("Z7").Value = $null
$ws.Range("AB7").Value = $null
